$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" query (row 2, column B) is being corrected: the final
# RETURN line that surfaced a `Cohort` column (via an OPTIONAL MATCH on
# (co:cohort)) is removed, since cohort is no longer part of this report.
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE s.clinical_study_designation IN ['COTC007B']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesQuery

# Window chrome left by the editing session: bigger zoom, selection
# resting on the Samples row.
$excel.ActiveWindow.Zoom = 115
$ws.Range("B3").Select()
